$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B97").Value = 57
$ws.Range("C97").Value = 529.2619502544403

$ws.Range("B98").Value = 58
$ws.Range("C98").Value = 563.6283569335938

$ws.Range("B99").Value = 61
$ws.Range("C99").Value = 1635.281559705734

$ws.Range("B100").Value = 62
$ws.Range("C100").Value = 1070.729019641876
